# Daily_Motivations.xlsx — "Updating the daily Scores"
#
# 1) For every existing "sleep" row (Motivation_typ = "sleep"), set the
#    JKL (col N) and OS (col O) scores to TRUE (they were all FALSE).
# 2) For the 2025-02-18 "activity" row (row 54), set HealthQuest (col L)
#    to TRUE (it was FALSE).
# 3) Append a new day, 2025-02-20, with its three rows (sleep, activity,
#    weekly_activity) of per-person scores.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Fix JKL / OS columns on every "sleep" row ------------------------
$sleepRows = @(2, 5, 8, 11, 14, 17, 20, 23, 26, 29, 32, 35, 38, 41, 44, 47, 50, 53, 56)
foreach ($r in $sleepRows) {
    $ws.Cells.Item($r, 14).Value = $true   # column N = JKL
    $ws.Cells.Item($r, 15).Value = $true   # column O = OS
}

# --- 2) Fix HealthQuest on the 2025-02-18 activity row -------------------
$ws.Cells.Item(54, 12).Value = $true       # column L = HealthQuest

# --- 3) Append the new 2025-02-20 rows (59, 60, 61) -----------------------
function Set-DateCell($row) {
    # Leading apostrophe forces the text "2025-02-20" to stay a string
    # instead of Excel auto-converting it to a date value/serial, matching
    # how every other Date cell in this sheet is stored (a plain string).
    # ClearFormats() then drops the "quote prefix" marker the apostrophe
    # leaves behind, so the cell ends up styled exactly like its neighbors
    # (no explicit style index, General format).
    $ws.Cells.Item($row, 1).Value = "'2025-02-20"
    $ws.Cells.Item($row, 1).ClearFormats()
}

# Row 59: sleep
Set-DateCell 59
$ws.Cells.Item(59, 2).Value = "sleep"
$ws.Cells.Item(59, 3).Value = $true     # Sportfinke
$ws.Cells.Item(59, 4).Value = $false    # Taylor Atwood
$ws.Cells.Item(59, 5).Value = $true     # Summerbody25
$ws.Cells.Item(59, 6).Value = $true     # Iron Man
$ws.Cells.Item(59, 7).Value = $true     # GurkenSalat
$ws.Cells.Item(59, 8).Value = $true     # yKing
$ws.Cells.Item(59, 9).Value = $true     # StayStrong
$ws.Cells.Item(59, 10).Value = $true    # WobblyWheel
$ws.Cells.Item(59, 11).Value = $true    # ClearMind23
$ws.Cells.Item(59, 12).Value = $false   # HealthQuest
$ws.Cells.Item(59, 13).Value = $true    # DeadliftCarror
$ws.Cells.Item(59, 14).Value = $true    # JKL
$ws.Cells.Item(59, 15).Value = $true    # OS

# Row 60: activity
Set-DateCell 60
$ws.Cells.Item(60, 2).Value = "activity"
$ws.Cells.Item(60, 3).Value = $false    # Sportfinke
$ws.Cells.Item(60, 4).Value = $false    # Taylor Atwood
$ws.Cells.Item(60, 5).Value = $false    # Summerbody25
$ws.Cells.Item(60, 6).Value = $false    # Iron Man
$ws.Cells.Item(60, 7).Value = $false    # GurkenSalat
$ws.Cells.Item(60, 8).Value = $true     # yKing
$ws.Cells.Item(60, 9).Value = $true     # StayStrong
$ws.Cells.Item(60, 10).Value = $true    # WobblyWheel
$ws.Cells.Item(60, 11).Value = $false   # ClearMind23
$ws.Cells.Item(60, 12).Value = $true    # HealthQuest
$ws.Cells.Item(60, 13).Value = $true    # DeadliftCarror
$ws.Cells.Item(60, 14).Value = $false   # JKL
$ws.Cells.Item(60, 15).Value = $false   # OS

# Row 61: weekly_activity
Set-DateCell 61
$ws.Cells.Item(61, 2).Value = "weekly_activity"
$ws.Cells.Item(61, 3).Value = $false    # Sportfinke
$ws.Cells.Item(61, 4).Value = $false    # Taylor Atwood
$ws.Cells.Item(61, 5).Value = $true     # Summerbody25
$ws.Cells.Item(61, 6).Value = $false    # Iron Man
$ws.Cells.Item(61, 7).Value = $true     # GurkenSalat
$ws.Cells.Item(61, 8).Value = $false    # yKing
$ws.Cells.Item(61, 9).Value = $true     # StayStrong
$ws.Cells.Item(61, 10).Value = $true    # WobblyWheel
$ws.Cells.Item(61, 11).Value = $false   # ClearMind23
$ws.Cells.Item(61, 12).Value = $false   # HealthQuest
$ws.Cells.Item(61, 13).Value = $true    # DeadliftCarror
$ws.Cells.Item(61, 14).Value = $false   # JKL
$ws.Cells.Item(61, 15).Value = $false   # OS
